$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1995073891625616
$ws.Range("C2").Value = 0.5517241379310345
$ws.Range("J2").Value = 0.004926108374384237
$ws.Range("P2").Value = 0.1428571428571428
$ws.Range("S2").Value = 0.1009852216748768
$ws.Range("C3").Value = 0.004504504504504504
$ws.Range("J3").Value = 0.02702702702702703
$ws.Range("P3").Value = 0.7882882882882883
$ws.Range("S3").Value = 0.1801801801801802
$ws.Range("B6").Value = 0.07048458149779736
$ws.Range("D6").Value = 0.013215859030837
$ws.Range("E6").Value = 0.004405286343612335
$ws.Range("F6").Value = 0.04405286343612335
$ws.Range("J6").Value = 0.2555066079295154
$ws.Range("O6").Value = 0.00881057268722467
$ws.Range("Q6").Value = 0.13215859030837
$ws.Range("R6").Value = 0.1145374449339207
$ws.Range("S6").Value = 0.3568281938325991
$ws.Range("B7").Value = 0.1268292682926829
$ws.Range("D7").Value = 0.01463414634146342
$ws.Range("F7").Value = 0.02926829268292683
$ws.Range("J7").Value = 0.1317073170731707
$ws.Range("O7").Value = 0.01463414634146342
$ws.Range("Q7").Value = 0.1414634146341463
$ws.Range("R7").Value = 0.07317073170731707
$ws.Range("S7").Value = 0.4682926829268293
$ws.Range("B8").Value = 0.1356673960612692
$ws.Range("D8").Value = 0.0175054704595186
$ws.Range("E8").Value = 0.006564551422319475
$ws.Range("F8").Value = 0.05032822757111598
$ws.Range("J8").Value = 0.1115973741794311
$ws.Range("O8").Value = 0.01969365426695843
$ws.Range("Q8").Value = 0.09190371991247265
$ws.Range("R8").Value = 0.09846827133479212
$ws.Range("S8").Value = 0.4682713347921225
$ws.Range("B9").Value = 0.1173469387755102
$ws.Range("D9").Value = 0.01020408163265306
$ws.Range("F9").Value = 0.04591836734693878
$ws.Range("J9").Value = 0.1071428571428571
$ws.Range("O9").Value = 0.01530612244897959
$ws.Range("Q9").Value = 0.1377551020408163
$ws.Range("R9").Value = 0.09693877551020408
$ws.Range("S9").Value = 0.4693877551020408
$ws.Range("B10").Value = 0.1289456010745467
$ws.Range("D10").Value = 0.02484889187374077
$ws.Range("F10").Value = 0.07387508394895903
$ws.Range("J10").Value = 0.1175285426460712
$ws.Range("O10").Value = 0.01477501678979181
$ws.Range("Q10").Value = 0.1746138347884486
$ws.Range("R10").Value = 0.07790463398253862
$ws.Range("S10").Value = 0.3875083948959033
$ws.Range("G11").Value = 0.1314102564102564
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.1858974358974359
$ws.Range("L11").Value = 0.5897435897435898
$ws.Range("S11").Value = 0.01602564102564102
$ws.Range("G12").Value = 0.6947368421052632
$ws.Range("J12").Value = 0.2263157894736842
$ws.Range("L12").Value = 0.02105263157894737
$ws.Range("S12").Value = 0.05789473684210526
$ws.Range("G13").Value = 0.7659574468085106
$ws.Range("J13").Value = 0.2127659574468085
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.0228310502283105
$ws.Range("H15").Value = 0.0958904109589041
$ws.Range("I15").Value = 0.091324200913242
$ws.Range("J15").Value = 0.3607305936073059
$ws.Range("K15").Value = 0.0730593607305936
$ws.Range("O15").Value = 0.0502283105022831
$ws.Range("S15").Value = 0.3059360730593607
$ws.Range("F16").Value = 0.01915708812260536
$ws.Range("H16").Value = 0.157088122605364
$ws.Range("I16").Value = 0.08812260536398467
$ws.Range("J16").Value = 0.4367816091954023
$ws.Range("K16").Value = 0.08045977011494253
$ws.Range("M16").Value = 0.01532567049808429
$ws.Range("O16").Value = 0.03448275862068965
$ws.Range("S16").Value = 0.1685823754789272
$ws.Range("F17").Value = 0.01799485861182519
$ws.Range("H17").Value = 0.1439588688946015
$ws.Range("I17").Value = 0.06169665809768637
$ws.Range("J17").Value = 0.5141388174807198
$ws.Range("K17").Value = 0.08226221079691516
$ws.Range("M17").Value = 0.0102827763496144
$ws.Range("O17").Value = 0.04884318766066838
$ws.Range("S17").Value = 0.1208226221079692
$ws.Range("F18").Value = 0.01363636363636364
$ws.Range("H18").Value = 0.1863636363636364
$ws.Range("I18").Value = 0.09545454545454546
$ws.Range("J18").Value = 0.3772727272727273
$ws.Range("K18").Value = 0.07727272727272727
$ws.Range("M18").Value = 0.02727272727272727
$ws.Range("O18").Value = 0.03636363636363636
$ws.Range("S18").Value = 0.1863636363636364
$ws.Range("F19").Value = 0.01519949335022166
$ws.Range("H19").Value = 0.1918936035465484
$ws.Range("I19").Value = 0.06966434452184928
$ws.Range("J19").Value = 0.3818872704243192
$ws.Range("K19").Value = 0.1000633312222926
$ws.Range("M19").Value = 0.02343255224825839
$ws.Range("N19").Value = 0.001266624445851805
$ws.Range("O19").Value = 0.06586447118429385
$ws.Range("S19").Value = 0.1507283090563648
